# Update the localization-status report with a fresh handoff run:
# the source file's generated GUID name changed from
# 413e9ea3-6742-4631-aa99-216f47bef52f to c1695c74-d002-440d-990b-85ecfacc4b2a,
# and the handoff/generate timestamps advanced a bit.

$wb = $excel.ActiveWorkbook

$oldGuid = "413e9ea3-6742-4631-aa99-216f47bef52f"
$newGuid = "c1695c74-d002-440d-990b-85ecfacc4b2a"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# The hyperlink's underlying target (relationship address) is left pointing
# at the original file/commit on GitHub; only the on-sheet display text is
# refreshed to the new generated file name.
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4293780bfee61738f45d215c462a1eaac3fed259/e2e/$oldGuid.md"

# --- Overview sheet ---
$ws1.Range("A2").Value = "$newGuid.md"

# B2 carries both the cell text and a hyperlink pointing at the file on GitHub;
# rebuild the hyperlink so its display text reflects the new file name.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $ghBase, $null, $null, "e2e\$newGuid.md")

# Latest HO Xliff Generate Date
$ws1.Range("G2").Value = "2016-09-06 05:14:11"
$ws1.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# --- zh-cn sheet ---
$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ghBase, $null, $null, "$newGuid.md")
$ws2.Range("G2").Value = "$newGuid.b062981febd2d0c6311f4aa71852a5b0a89c1bdf.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-06 05:14:01"
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# --- de-de sheet ---
$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ghBase, $null, $null, "$newGuid.md")
$ws3.Range("G2").Value = "$newGuid.b062981febd2d0c6311f4aa71852a5b0a89c1bdf.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-06 05:14:11"
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
